$wb = $excel.ActiveWorkbook

# Reference cells already carrying the workbook's standard "header" /
# "row-index" cell format (bold, centered, thin border) so the new sheets
# end up stylistically identical to the existing quarterly sheets.
$styleSrc = $wb.Worksheets.Item("2021-Q4")
$headerStyleCell = $styleSrc.Range("B1")
$indexStyleCell = $styleSrc.Range("A2")

function Copy-Format($srcCell, $destRange) {
    $srcCell.Copy()
    $destRange.PasteSpecial(-4122)  # xlPasteFormats
}

# Helper: write a value as TEXT (the source data stores numbers such as
# "10.92" / "519029" as text, not numbers) without leaving a stray style
# behind because of the quote-prefix trick used to force text.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# --- Step 1: the current "总计" sheet becomes the new "2022-Q1" sheet -------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Wipe the old "总计" data/formatting out of this sheet before writing the
# new quarterly fund-holdings table into it.
$q1.Cells.Clear()

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"
Copy-Format $headerStyleCell $q1.Range("B1:H1")

# Row 2
$q1.Range("A2").Value = 0
Set-TextValue $q1.Range("B2") "519029"
$q1.Range("C2").Value = "华夏稳增混合"
Set-TextValue $q1.Range("D2") "10.92"
Set-TextValue $q1.Range("E2") "92.99"
Set-TextValue $q1.Range("F2") "6.79"
Set-TextValue $q1.Range("G2") "0.7415"
$q1.Range("H2").Value = 1
Copy-Format $indexStyleCell $q1.Range("A2")

# Row 3
$q1.Range("A3").Value = 1
Set-TextValue $q1.Range("B3") "003501"
$q1.Range("C3").Value = "泰达宏利睿智稳健灵活配置混合"
Set-TextValue $q1.Range("D3") "3.57"
Set-TextValue $q1.Range("E3") "73.79"
Set-TextValue $q1.Range("F3") "2.08"
Set-TextValue $q1.Range("G3") "0.0743"
$q1.Range("H3").Value = 9
Copy-Format $indexStyleCell $q1.Range("A3")

# --- Step 2: append a brand-new "总计" sheet after the last sheet ----------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$total = $wb.Worksheets.Add($null, $lastSheet)
$total.Name = "总计"

# Header row
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"
Copy-Format $headerStyleCell $total.Range("B1:D1")

# Data rows (newest quarter first, same ordering the sheet already used)
$rows = @(
    @(0, "2022-Q1", 2, 0.82),
    @(1, "2021-Q4", 2, 0.79),
    @(2, "2021-Q3", 1, 0.09),
    @(3, "2021-Q1", 3, 0.99),
    @(4, "2020-Q4", 2, 0.23)
)

$r = 2
foreach ($row in $rows) {
    $total.Range("A$r").Value = $row[0]
    $total.Range("B$r").Value = $row[1]
    $total.Range("C$r").Value = $row[2]
    $total.Range("D$r").Value = $row[3]
    Copy-Format $indexStyleCell $total.Range("A$r")
    $r = $r + 1
}

# Restore the originally-active sheet/tab (unchanged by this edit).
$wb.Worksheets.Item(1).Activate()
